# Applies newly-collected HTTPS local client-side-time streaming
# measurement readings (100K / 1000K rows) to sheets
# "1 StreamingSources" .. "5 StreamingSources".
#
# Previously these rows were blank, which made the downstream AVERAGE()
# formulas evaluate to #DIV/0! (propagating into the "summary" block at
# the bottom of each sheet, and showing as 0 in the charts). Filling in
# the readings lets everything recalc normally.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "1 StreamingSources" — columns B (Time secs), C (Memory Mbs)
#   row 18-20 -> "1000K" group, row 22-24 -> "2000K" group
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("1 StreamingSources")
$ws1.Range("B18").Value = 2111.0430000000001
$ws1.Range("C18").Value = 123
$ws1.Range("B19").Value = 2042.174
$ws1.Range("C19").Value = 79
$ws1.Range("B22").Value = 4084.0210000000002
$ws1.Range("C22").Value = 97

# ---------------------------------------------------------------------
# "2 StreamingSources" — columns B,C (Time secs), D (Memory Mbs)
#   row 14-16 -> "100K" group, row 18 -> "1000K" group
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2 StreamingSources")
$ws2.Range("B14").Value = 208.40299999999999
$ws2.Range("C14").Value = 207.74799999999999
$ws2.Range("D14").Value = 53
$ws2.Range("B15").Value = 205.72800000000001
$ws2.Range("C15").Value = 205.68600000000001
$ws2.Range("D15").Value = 127
$ws2.Range("B16").Value = 205.29499999999999
$ws2.Range("C16").Value = 205.10300000000001
$ws2.Range("D16").Value = 118
$ws2.Range("B18").Value = 2045.24
$ws2.Range("C18").Value = 2045.26
$ws2.Range("D18").Value = 121

# ---------------------------------------------------------------------
# "3 StreamingSources" — columns B,C,D (Time secs), E (Memory Mbs)
#   row 14-16 -> "100K" group, row 18 -> "1000K" group
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("3 StreamingSources")
$ws3.Range("B14").Value = 208.62799999999999
$ws3.Range("C14").Value = 208.56399999999999
$ws3.Range("D14").Value = 209.11699999999999
$ws3.Range("E14").Value = 130
$ws3.Range("B15").Value = 207.09100000000001
$ws3.Range("C15").Value = 206.946
$ws3.Range("D15").Value = 206.59899999999999
$ws3.Range("E15").Value = 126
$ws3.Range("B16").Value = 207.07300000000001
$ws3.Range("C16").Value = 206.523
$ws3.Range("D16").Value = 207.34800000000001
$ws3.Range("E16").Value = 130
$ws3.Range("B18").Value = 2083.7939999999999
$ws3.Range("C18").Value = 2075.502
$ws3.Range("D18").Value = 2075.364
$ws3.Range("E18").Value = 161

# ---------------------------------------------------------------------
# "4 StreamingSources" — columns B,C,D,E (Time secs), F (Memory Mbs)
#   row 14-16 -> "100K" group, row 18 -> "1000K" group
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("4 StreamingSources")
$ws4.Range("B14").Value = 250.47300000000001
$ws4.Range("C14").Value = 249.76400000000001
$ws4.Range("D14").Value = 246.51
$ws4.Range("E14").Value = 249.78200000000001
$ws4.Range("F14").Value = 140
$ws4.Range("B15").Value = 220.71899999999999
$ws4.Range("C15").Value = 216.35900000000001
$ws4.Range("D15").Value = 218.10499999999999
$ws4.Range("E15").Value = 217.96700000000001
$ws4.Range("F15").Value = 137
$ws4.Range("B16").Value = 218.02099999999999
$ws4.Range("C16").Value = 213.77
$ws4.Range("D16").Value = 216.37299999999999
$ws4.Range("E16").Value = 215.81399999999999
$ws4.Range("F16").Value = 146
$ws4.Range("B18").Value = 2167.1170000000002
$ws4.Range("C18").Value = 2152.9389999999999
$ws4.Range("D18").Value = 2175.0349999999999
$ws4.Range("E18").Value = 2187.4490000000001
$ws4.Range("F18").Value = 159

# ---------------------------------------------------------------------
# "5 StreamingSources" — columns B,C,D,E,F (Time secs), G (Memory Mbs)
#   row 14-16 -> "100K" group, row 18 -> "1000K" group
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("5 StreamingSources")
$ws5.Range("B14").Value = 287.58
$ws5.Range("C14").Value = 276.07100000000003
$ws5.Range("D14").Value = 262.00700000000001
$ws5.Range("E14").Value = 270.65100000000001
$ws5.Range("F14").Value = 279.56700000000001
$ws5.Range("G14").Value = 158
$ws5.Range("B15").Value = 270.89499999999998
$ws5.Range("C15").Value = 270.089
$ws5.Range("D15").Value = 269.96199999999999
$ws5.Range("E15").Value = 264.029
$ws5.Range("F15").Value = 270.86399999999998
$ws5.Range("G15").Value = 160
$ws5.Range("B16").Value = 278.00700000000001
$ws5.Range("C16").Value = 272.34399999999999
$ws5.Range("D16").Value = 277.33100000000002
$ws5.Range("E16").Value = 270.10300000000001
$ws5.Range("F16").Value = 278.298
$ws5.Range("G16").Value = 152
$ws5.Range("B18").Value = 2527.1019999999999
$ws5.Range("C18").Value = 2500.0239999999999
$ws5.Range("D18").Value = 2493.44
$ws5.Range("E18").Value = 2453.7460000000001
$ws5.Range("F18").Value = 2531.6280000000002
$ws5.Range("G18").Value = 224

# ---------------------------------------------------------------------
# Restore the final cursor/selection position on each affected sheet to
# match where the author was last working when the workbook was saved.
# ---------------------------------------------------------------------
$ws2.Range("B19").Select()
$ws3.Range("B19").Select()
$ws4.Range("B19").Select()
$ws5.Range("B19").Select()

$ws1.Activate()
$ws1.Range("Q28").Select()
